$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = @"
TC_LOGIN_02
"@
$ws.Range("B3").Value = @"
Login
"@
$ws.Range("C3").Value = @"
Login with valid email and invalid password
"@
$ws.Range("D3").Value = @"
1. Enter valid email
2. Enter invalid password
3. Click Login
"@
$ws.Range("E3").Value = @"
Email:superadmin
Password:admin
"@
$ws.Range("F3").Value = @"
Error message should be displayed
"@
$ws.Range("G3").Value = @"
Failed
"@
$ws.Range("H3").Value = @"
Expected: Error message should be displayed | Actual: Locator expected to be visible
Actual value: <element(s) not found> 
Call log:
  - Expect "to_be_visible" with timeout 5000ms
  - waiting for locator(".error, .error-message, .alert-danger, .validation-error")

"@

$ws.Range("A4").Value = @"
TC_LOGIN_03
"@
$ws.Range("B4").Value = @"
Login
"@
$ws.Range("C4").Value = @"
Login with invalid email and valid password
"@
$ws.Range("D4").Value = @"
1. Enter invalid email
2. Enter valid password
3. Click Login
"@
$ws.Range("E4").Value = @"
Email:super
Password:superadmin
"@
$ws.Range("F4").Value = @"
Error message should be displayed
"@
$ws.Range("G4").Value = @"
Failed
"@
$ws.Range("H4").Value = @"
Expected: Error message should be displayed | Actual: Locator expected to be visible
Actual value: <element(s) not found> 
Call log:
  - Expect "to_be_visible" with timeout 5000ms
  - waiting for locator(".error, .error-message, .alert-danger, .validation-error")

"@

$ws.Range("A5").Value = @"
TC_LOGIN_04
"@
$ws.Range("B5").Value = @"
Login
"@
$ws.Range("C5").Value = @"
Login with invalid email and invalid password
"@
$ws.Range("D5").Value = @"
1. Enter invalid email
2. Enter invalid password
3. Click Login
"@
$ws.Range("E5").Value = @"
Email:superadmin
Password:superadmin
"@
$ws.Range("F5").Value = @"
Error message should be displayed
"@
$ws.Range("G5").Value = @"
Failed
"@
$ws.Range("H5").Value = @"
Expected: Error message should be displayed | Actual: Locator expected to be visible
Actual value: <element(s) not found> 
Call log:
  - Expect "to_be_visible" with timeout 5000ms
  - waiting for locator(".error, .error-message, .alert-danger, .validation-error")

"@

$ws.Range("A6").Value = @"
TC_LOGIN_05
"@
$ws.Range("B6").Value = @"
Login
"@
$ws.Range("C6").Value = @"
Login with empty email and empty password
"@
$ws.Range("D6").Value = @"
1. Leave fields empty
2. Click Login
"@
$ws.Range("E6").Value = @"
Email:superadmin
Password:superadmin
"@
$ws.Range("F6").Value = @"
Validation message should be shown
"@
$ws.Range("G6").Value = @"
Failed
"@
$ws.Range("H6").Value = @"
Expected: Validation message should be shown | Actual: Locator expected to be visible
Actual value: <element(s) not found> 
Call log:
  - Expect "to_be_visible" with timeout 5000ms
  - waiting for locator(".error, .error-message, .alert-danger, .validation-error")

"@

$ws.Range("A7").Value = @"
TC_LOGIN_06
"@
$ws.Range("B7").Value = @"
Login
"@
$ws.Range("C7").Value = @"
Login with empty email
"@
$ws.Range("D7").Value = @"
1. Leave email empty
2. Enter password
3. Click Login
"@
$ws.Range("E7").Value = @"
Email:superadmin
Password:superadmin
"@
$ws.Range("F7").Value = @"
Email required validation should appear
"@
$ws.Range("G7").Value = @"
Failed
"@
$ws.Range("H7").Value = @"
Expected: Email required validation should appear | Actual: Locator expected to be visible
Actual value: <element(s) not found> 
Call log:
  - Expect "to_be_visible" with timeout 5000ms
  - waiting for locator(".error, .error-message, .alert-danger, .validation-error")

"@

$ws.Range("A8").Value = @"
TC_LOGIN_07
"@
$ws.Range("B8").Value = @"
Login
"@
$ws.Range("C8").Value = @"
Login with empty password
"@
$ws.Range("D8").Value = @"
1. Enter email
2. Leave password empty
3. Click Login
"@
$ws.Range("E8").Value = @"
Email:superadmin
Password:superadmin
"@
$ws.Range("F8").Value = @"
Password required validation should appear
"@
$ws.Range("G8").Value = @"
Failed
"@
$ws.Range("H8").Value = @"
Expected: Password required validation should appear | Actual: Locator expected to be visible
Actual value: <element(s) not found> 
Call log:
  - Expect "to_be_visible" with timeout 5000ms
  - waiting for locator(".error, .error-message, .alert-danger, .validation-error")

"@
